$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Action" notes added in column H, next to the role rows (System / Processor block).
$ws.Range("H4").Value = "Huy: design generate report (create, review); core module"
$ws.Range("H5").Value = "Thien: mobile app for truck driver, design UI, find CSS template, web services for mobile app"
$ws.Range("H6").Value = "Thinh: mobile app for goods owner, web services for mobile app"
$ws.Range("H7").Value = "Khuong: web app, webservice for web"

# Size column H to fit its (longest) new content (best-fit width).
$ws.Columns.Item(8).ColumnWidth = 84

# Leave the selection where the edits ended, matching the saved view state.
$ws.Range("H7").Select() | Out-Null
